# The 12.b.1 indicator description (cell B4 on the single worksheet) was
# rewritten to the officially revised wording. Updating the cell's text
# removes the old shared-string entry and appends the new one at the end
# of the shared-strings table (standard Excel behaviour for a changed
# literal string cell), which also re-numbers the shared-string indices
# used by the other (unchanged) cells below it - that ripple is expected
# and not something to special-case here.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "12.b.1 Внедрение стандартных методов учета в целях отслеживания экономических и экологических характеристик устойчивости туризма"

# Leave the active selection on the cell that was edited, matching the
# saved workbook's cursor position.
$ws.Range("B4").Select()
